# Update capital structure database values for Sudan bank rows (2 and 3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("D$row").Value = 0.475
    $ws.Range("E$row").Value = 0.599
    $ws.Range("K$row").Value = 7.41
    $ws.Range("L$row").Value = 0.5292857142857142
    $ws.Range("U$row").Value = 32.8
    $ws.Range("V$row").Value = 0.6096654275092936
    $ws.Range("W$row").Value = 0.1193236714975845
    $ws.Range("X$row").Value = 0.1238467486970918
    $ws.Range("Y$row").Value = -0.004523077199507286
    $ws.Range("Z$row").Value = 0.392156862745098
    $ws.Range("AB$row").Value = 0.1238467486970918
    $ws.Range("AC$row").Value = -0.1238467486970918
    $ws.Range("AG$row").Value = -32.8
    $ws.Range("AJ$row").Value = -1.561904761904762
    $ws.Range("AK$row").Value = -1.163120567375886
}
